$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new abbreviation rows (61-71) ---
# Column "Initials" (A) is filled first for rows 61-70, then column "Complete
# words" (B) for rows 61-65, then the last initial (A71), then the remaining
# complete words (B66-B71) -- this mirrors the original authoring order so
# the shared-string table comes out in the same sequence.

$ws.Range("A61").Value = "HVM"
$ws.Range("A62").Value = "HVP"
$ws.Range("A63").Value = "HTHP"
$ws.Range("A64").Value = "HPHT"
$ws.Range("A65").Value = "PPT"
$ws.Range("A66").Value = "PSD"
$ws.Range("A67").Value = "PPB"
$ws.Range("A68").Value = "BBL"
$ws.Range("A69").Value = "CONC"
$ws.Range("A70").Value = "DIR"

$ws.Range("B61").Value = "High viscosity Mud"
$ws.Range("B62").Value = "High viscosity pill"
$ws.Range("B63").Value = "Hight temperature High pressure"
$ws.Range("B64").Value = "High pressure High temperature"
$ws.Range("B65").Value = "Plugging permeability test"

$ws.Range("A71").Value = "PPA"

$ws.Range("B66").Value = "particle size distribution"
$ws.Range("B67").Value = "pound per barrel"
$ws.Range("B68").Value = "barrels"
$ws.Range("B69").Value = "Concentration"
$ws.Range("B70").Value = "Directional"
$ws.Range("B71").Value = "Plugging permeability apparatus"

# --- Grow the table / autofilter range to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B71"))

# --- Widen column B to fit the longer phrases ---
$ws.Columns.Item(2).ColumnWidth = 31.6666666666667

# --- Leave the selection where the author ended up ---
$ws.Range("B72").Select() | Out-Null
